# Swap columns L (Se) and M (As), and columns R (Ag) and S (Pb),
# for the header row and all data rows (1-16), effectively reordering
# the element columns "As" before "Se" and "Pb" before "Ag".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 1
$lastRow = $ws.UsedRange.Rows.Count

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # Swap column L and M
    $lVal = $ws.Cells.Item($r, 12).Value2
    $mVal = $ws.Cells.Item($r, 13).Value2
    $ws.Cells.Item($r, 12).Value2 = $mVal
    $ws.Cells.Item($r, 13).Value2 = $lVal

    # Swap column R and S
    $rVal = $ws.Cells.Item($r, 18).Value2
    $sVal = $ws.Cells.Item($r, 19).Value2
    $ws.Cells.Item($r, 18).Value2 = $sVal
    $ws.Cells.Item($r, 19).Value2 = $rVal
}
